$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-24 Saturday" "2024-02-25 Sunday"

Replace-Text "62÷2=" "88÷4="
Replace-Text "30÷5=" "94÷9="
Replace-Text "33÷5=" "48÷6="
Replace-Text "80÷3=" "64÷5="
Replace-Text "66÷8=" "43÷4="
Replace-Text "81÷8=" "63÷4="
Replace-Text "23÷2=" "74÷6="
Replace-Text "89÷3=" "87÷7="
Replace-Text "72÷8=" "66÷7="
Replace-Text "48÷5=" "98÷2="
Replace-Text "94÷5=" "41÷6="
Replace-Text "12÷8=" "53÷3="
Replace-Text "91÷3=" "31÷5="
Replace-Text "90÷2=" "87÷5="
Replace-Text "71÷9=" "50÷9="
Replace-Text "16÷8=" "40÷5="
Replace-Text "52÷4=" "27÷6="
Replace-Text "10÷9=" "34÷7="
Replace-Text "37÷8=" "43÷3="
Replace-Text "17÷2=" "68÷2="
Replace-Text "26÷9=" "87÷8="
Replace-Text "45÷7=" "71÷7="
Replace-Text "20÷2=" "12÷6="
Replace-Text "32÷3=" "35÷5="
Replace-Text "52÷2=" "76÷4="
